# edit.ps1
#
# Commit: "edit simple graph slides"
#
# The authoritative diff shows two adjacent runs inside one paragraph being
# collapsed into a single run:
#
#   <a:r><a:rPr .../><a:t> </a:t></a:r><a:r><a:rPr .../><a:t>8,</a:t></a:r>
#   -->
#   <a:r><a:rPr .../><a:t> 8,</a:t></a:r>
#
# Both runs already shared identical run properties (sz="1200" baseline="0"
# Comic Sans MS), so the edit is simply "retype"/merge of the space and the
# "8," that follows it into one run - the kind of no-visible-effect cleanup
# PowerPoint itself performs when you edit text that happens to span a run
# boundary with identical formatting.
#
# This script reproduces that behaviour generically: it scans every
# slide/shape (including nested group items and table cells) for a run of
# text " 8," where the leading space and the "8," are adjacent characters,
# then rewrites that character span as a single piece of text so the
# underlying runs collapse into one - exactly mirroring the diff - while
# leaving every other run/shape untouched.

function Get-CharKey($chars) {
    $f = $chars.Font
    $rgb = $null
    try { $rgb = $f.Color.RGB } catch {}
    return ($f.Name + "|" + $f.Size + "|" + $f.Bold + "|" + $f.Italic + "|" + $f.Underline + "|" + [string]$rgb)
}

function Fix-TextRange($tr) {
    if ($tr -eq $null) { return }
    $text = $tr.Text
    if ($text -eq $null) { return }
    if ($text.Length -lt 3) { return }

    # Look for every occurrence of "8," immediately preceded by a run whose
    # entire text is a single space (matching the diff: a lone " " run
    # directly followed by an "8," run, both sharing identical formatting).
    # Merging that 3-character span into one literal " 8," run collapses the
    # two original runs into a single run, exactly like the diff.
    $searchFrom = 1
    while ($true) {
        $idx = $text.IndexOf("8,", $searchFrom - 1)
        if ($idx -lt 0) { break }
        $pos = $idx + 1   # 1-based position of "8"

        if ($pos -gt 1) {
            $prevChar = $tr.Characters($pos - 1, 1)
            if ($prevChar.Text -eq " ") {
                $digitRun = $tr.Characters($pos, 2)
                $prevKey = Get-CharKey $prevChar
                $digitKey = Get-CharKey $digitRun

                # Make sure the space is its own run (i.e. not merely the
                # tail character of a longer run) by checking that the
                # character before it - if any - does NOT share the same
                # formatting; a real run boundary separates them.
                $isLoneSpaceRun = $true
                if ($pos -gt 2) {
                    $beforeSpace = $tr.Characters($pos - 2, 1)
                    $beforeKey = Get-CharKey $beforeSpace
                    if ($beforeKey -eq $prevKey) {
                        $isLoneSpaceRun = $false
                    }
                }

                if ($isLoneSpaceRun -and ($prevKey -eq $digitKey)) {
                    $merged = $tr.Characters($pos - 1, 3)
                    if ($merged.Text -eq " 8,") {
                        $merged.Text = " 8,"
                    }
                }
            }
        }

        $searchFrom = $idx + 2
        $text = $tr.Text
        if ($searchFrom -gt $text.Length) { break }
    }
}

function Walk-Shape($shape) {
    if ($shape -eq $null) { return }

    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        Fix-TextRange $shape.TextFrame.TextRange
    }

    if ($shape.Type -eq 6 -or $shape.GroupItems.Count -gt 0) {
        # msoGroup == 6; guarded by GroupItems.Count in case Type isn't exposed
        for ($i = 1; $i -le $shape.GroupItems.Count; $i++) {
            Walk-Shape $shape.GroupItems.Item($i)
        }
    }
}

$p = $ppt.ActivePresentation
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        Walk-Shape $slide.Shapes.Item($shi)
    }
}
